# psychopy loads its stimulus images relative to the experiment folder, so
# the "image" column in the conditions table needs to point into the
# "images" subfolder rather than to bare filenames.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$used = $ws.UsedRange
for ($r = 1; $r -le $used.Rows.Count; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -eq "up.jpg") {
        $cell.Value = "images/up.jpg"
    } elseif ($val -eq "down.jpg") {
        $cell.Value = "images/down.jpg"
    }
}

# The longer path text no longer fits the old row height for the row that
# was just edited, so its row grows to match the rest of the data rows.
$ws.Rows.Item(2).RowHeight = $ws.Rows.Item(3).RowHeight

# The "image" column now holds longer "images/....jpg" paths, so it needs
# to be widened to keep showing the full text (target stored width is
# ~15.21 chars; ColumnWidth excludes the standard ~5px/0.83-char grid
# padding that Excel adds on top when it writes the column definition).
$ws.Columns.Item(2).ColumnWidth = 15.2117647058824 - 0.8333333333333334

# Leave the cursor where the edit happened.
$ws.Activate()
$ws.Range("C17").Select() | Out-Null
